$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K" - strikeouts) values for rows 2-9
$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 4
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 3
$ws.Range("G9").Value = 3
